$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.580.93'
$ws.Range('E2').Value = '  -1.83%  '
$ws.Range('D3').Value = '3.830.61'
$ws.Range('E3').Value = '  -2.43%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.25'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.58'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').Value = '3.824.71'
$ws.Range('E7').Value = '  -2.51%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.523'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.67%  '
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.25'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.455'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000246'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.70'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.69%  '
$ws.Range('D15').Value = '4.471.81'
$ws.Range('E15').Value = '  -2.39%  '
$ws.Range('D16').Value = '3.835.74'
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('D17').Value = '67.678.11'
$ws.Range('E17').Value = '  -1.89%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.42'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.06'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +5.78%  '
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.68'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '466.65'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -4.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.726'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000158'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.22'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.19'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.06'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.94'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.25%  '
$ws.Range('E30').Value = '  -1.80%  '
$ws.Range('D31').Value = '3.978.65'
$ws.Range('E31').Value = '  -2.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.66'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.28'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -5.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '30.79'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.77%  '
$ws.Range('D35').Value = '3.801.14'
$ws.Range('E35').Value = '  -1.87%  '
$ws.Range('E36').Value = '  -3.52%  '
$ws.Range('E37').Value = '  -0.99%  '
$ws.Range('E38').Value = '  -3.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.86'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.24'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +7.68%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.310'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.78%  '
$ws.Range('E43').Value = '  -1.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '421.67'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '47.16'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.51'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '142.86'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.77%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000265'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +12.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0354'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '38.93'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.20%  '
